# Applies the "Penalty Reward System" forecast-refresh edit:
#  - Forecast Comparison: each week's Week_Start_Date rolls forward by one
#    week (row 2's old W2 date becomes row 2's date, etc.) and MyForecast
#    (column D) gets refreshed numbers.
#  - Summary: recomputed aggregate stats that follow from the above.
#
# Date-shaped and purely-numeric-shaped strings are written via
# Range.Formula with a leading apostrophe so Excel stores them as literal
# text (matching the source workbook's inline-string cells) instead of
# silently parsing them into date serials / numbers.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison: Week_Start_Date (B) and MyForecast (D) ---

$ws1.Range("B2").Formula = "'2025-01-12"
$ws1.Range("D2").Value = 32

$ws1.Range("B3").Formula = "'2025-01-19"
$ws1.Range("D3").Value = 31

$ws1.Range("B4").Formula = "'2025-01-26"
$ws1.Range("D4").Value = 30

$ws1.Range("B5").Formula = "'2025-02-02"
$ws1.Range("D5").Value = 28

$ws1.Range("B6").Formula = "'2025-02-09"
$ws1.Range("D6").Value = 28

$ws1.Range("B7").Formula = "'2025-02-16"
$ws1.Range("D7").Value = 31

$ws1.Range("B8").Formula = "'2025-02-23"
$ws1.Range("D8").Value = 37

$ws1.Range("B9").Formula = "'2025-03-02"
$ws1.Range("D9").Value = 42

$ws1.Range("B10").Formula = "'2025-03-09"
$ws1.Range("D10").Value = 29

$ws1.Range("B11").Formula = "'2025-03-16"
$ws1.Range("D11").Value = 28

$ws1.Range("B12").Formula = "'2025-03-23"
$ws1.Range("D12").Value = 27

$ws1.Range("B13").Formula = "'2025-03-30"
$ws1.Range("D13").Value = 26

$ws1.Range("B14").Formula = "'2025-04-06"
$ws1.Range("D14").Value = 27

$ws1.Range("B15").Formula = "'2025-04-13"
$ws1.Range("D15").Value = 34

$ws1.Range("B16").Formula = "'2025-04-20"
$ws1.Range("D16").Value = 30

$ws1.Range("B17").Formula = "'2025-04-27"
$ws1.Range("D17").Value = 30

# --- Summary: recomputed aggregate metrics ---

$ws2.Range("B2").Value = "2023-01-01 to 2025-01-05"

$ws2.Range("B5").Formula = "'28"
$ws2.Range("B6").Formula = "'17"

$ws2.Range("B8").Value = "3146 units"

$ws2.Range("B9").Formula = "'491"
$ws2.Range("B10").Formula = "'259"
$ws2.Range("B11").Formula = "'121"

$ws2.Range("B13").Formula = "'2025-03-02"
$ws2.Range("B15").Formula = "'2025-03-30"
